$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Months (column A) repeating pattern already present in the sheet: Enero..Junio for the new year 2023
$months = @("Enero", "Febrero", "Marzo", "Abril", "Mayo", "Junio")
$year = 2023
$startRow = 50

for ($i = 0; $i -lt $months.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $months[$i]
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $ws.Cells.Item($row, 5).Value = $year
}

# Update the view to reflect where the user had scrolled/selected after the edit
$ws.Range("H48").Select()
$excel.ActiveWindow.ScrollRow = 16

# Update column B width to fit updated content (closest reachable value to 30.51 chars)
$ws.Columns.Item(2).ColumnWidth = 29.65
